# This script fixes BLACKBERRY_LTD.xlsx rows that had been accidentally
# populated with price/share data belonging to other tickers (an "extra
# files" mix-up). For each data row (2-43) the open/close/high/low price
# columns (D:G), the shares_outstanding column (H), and the fixed_ticker
# column (I) are corrected back to BlackBerry ("BB") values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
    ,(2, 10.94999980926514, 8.930000305175781, 11.44999980926514, 8.590000152587891, 590361084)
    ,(3, 9.819999694824221, 8.180000305175781, 9.939999580383301, 8.039999961853027, 590361084)
    ,(4, 7.460000038146973, 6.130000114440918, 7.820000171661377, 5.960000038146973, 590361084)
    ,(5, 7.960000038146973, 9.279999732971191, 9.460000038146973, 7.28000020980835, 590361084)
    ,(6, 7.900000095367432, 8.090000152587891, 8.359999656677246, 7.429999828338623, 590361084)
    ,(7, 7.21999979019165, 6.710000038146973, 7.489999771118164, 6.230000019073486, 590361084)
    ,(8, 7.590000152587891, 7.980000019073486, 8.460000038146973, 7.150000095367432, 590361084)
    ,(9, 7.75, 6.889999866485596, 8.050000190734863, 6.860000133514404, 590361084)
    ,(10, 7, 7.75, 8.079999923706055, 6.650000095367432, 590361084)
    ,(11, 10.65999984741211, 9.989999771118164, 11.73999977111816, 9.520000457763672, 590361084)
    ,(12, 9.329999923706056, 11.18000030517578, 11.34000015258789, 8.890000343322754, 590361084)
    ,(13, 10.85000038146973, 11.17000007629394, 12.35999965667725, 10.13000011444092, 590361084)
    ,(14, 12.14999961853027, 11.5, 13.38000011444092, 11.47999954223633, 590361084)
    ,(15, 11.9399995803833, 9.649999618530272, 12.55000019073486, 9.520000457763672, 590361084)
    ,(16, 10.55000019073486, 11.38000011444092, 12, 9.670000076293944, 590361084)
    ,(17, 8.970000267028809, 7.110000133514404, 9, 6.570000171661377, 590361084)
    ,(18, 8.710000038146973, 10.09000015258789, 10.28999996185303, 8.569999694824219, 590361084)
    ,(19, 7.829999923706055, 7.460000038146973, 9.090000152587891, 7.099999904632568, 590361084)
    ,(20, 6.789999961853027, 5.25, 7.760000228881836, 5.099999904632568, 590361084)
    ,(21, 5.570000171661377, 6.420000076293945, 6.650000095367432, 5.329999923706055, 590361084)
    ,(22, 5.170000076293945, 4.130000114440918, 5.21999979019165, 2.700000047683716, 590361084)
    ,(23, 4.639999866485596, 4.889999866485596, 5.840000152587891, 4.590000152587891, 590361084)
    ,(24, 5.25, 4.590000152587891, 5.659999847412109, 4.489999771118164, 590361084)
    ,(25, 8.359999656677246, 6.630000114440918, 9.689999580383301, 6.599999904632568, 590361084)
    ,(26, 10.28999996185303, 8.430000305175781, 12.52000045776367, 8.359999656677246, 590361084)
    ,(27, 11.22999954223633, 12.22000026702881, 20.17000007629395, 10.52000045776367, 590361084)
    ,(28, 11.44999980926514, 9.729999542236328, 12.01000022888184, 9.239999771118164, 590361084)
    ,(29, 9.689999580383301, 9.350000381469728, 9.689999580383301, 8.420000076293945, 590361084)
    ,(30, 6.880000114440918, 7.460000038146973, 7.96999979019165, 5.880000114440918, 590361084)
    ,(31, 6.039999961853027, 5.389999866485596, 6.210000038146973, 5, 590361084)
    ,(32, 5.889999866485596, 4.699999809265137, 6.210000038146973, 4.679999828338623, 590361084)
    ,(33, 4.880000114440918, 3.259999990463257, 5.28000020980835, 3.170000076293945, 590361084)
    ,(34, 3.880000114440918, 4.559999942779541, 4.699999809265137, 3.430000066757202, 590361084)
    ,(35, 5.300000190734863, 5.53000020980835, 5.71999979019165, 4.710000038146973, 590361084)
    ,(36, 5.619999885559082, 4.710000038146973, 5.75, 4.489999771118164, 590361084)
    ,(37, 3.680000066757202, 3.539999961853028, 4.440000057220459, 3.480000019073486, 590361084)
    ,(38, 2.789999961853028, 2.759999990463257, 3.180000066757202, 2.460000038146973, 590361084)
    ,(39, 2.920000076293945, 2.480000019073486, 2.950000047683716, 2.200000047683716, 590361084)
    ,(40, 2.339999914169312, 2.619999885559082, 2.630000114440918, 2.220000028610229, 590361084)
    ,(41, 2.579999923706055, 3.779999971389771, 4.059999942779541, 2.509999990463257, 590361084)
    ,(42, 4.769999980926514, 3.769999980926514, 4.960000038146973, 3.660000085830689, 590361084)
    ,(43, 3.930000066757202, 4.579999923706055, 5.320000171661377, 3.869999885559082, 590361084)
)

foreach ($row in $data) {
    $r          = $row[0]
    $openPrice  = $row[1]
    $closePrice = $row[2]
    $highPrice  = $row[3]
    $lowPrice   = $row[4]
    $shares     = $row[5]

    $ws.Range("D$r").Value = $openPrice
    $ws.Range("E$r").Value = $closePrice
    $ws.Range("F$r").Value = $highPrice
    $ws.Range("G$r").Value = $lowPrice
    $ws.Range("H$r").Value = $shares
    $ws.Range("I$r").Value = "BB"
}
